# Auto-generated Excel COM-interop script
# Applies: (1) shared-string date update in confidentiality footer cell A59
#          (2) updated Weight (D) / Percent Change (E) values for rows 2-56

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected; unprotect so values can be written, matching
# how Excel itself would require unprotecting before editing cells.
$ws.Unprotect()

# --- Update the confidentiality footer date (2021-05-10 -> 2021-05-11) ---
$footerCell = $ws.Range("A59")
$footerCell.Value = $footerCell.Value2.Replace("2021-05-10", "2021-05-11")

# --- Update Weight (D) / Percent Change (E) figures ---
$ws.Range("D2").Value = 0.01547884321747375
$ws.Range("E2").Value = -0.01427021161150321
$ws.Range("D3").Value = 0.04962228903852196
$ws.Range("E3").Value = 0.01047488003410146
$ws.Range("D4").Value = 0.01441687157814867
$ws.Range("E4").Value = 0.001143547586683002
$ws.Range("D5").Value = 0.009935376145296751
$ws.Range("E5").Value = -0.006386975579210863
$ws.Range("D6").Value = 0.01551735290960862
$ws.Range("E6").Value = -0.007292802617229976
$ws.Range("D7").Value = 0.02012734877720137
$ws.Range("E7").Value = -0.004939347715551601
$ws.Range("D8").Value = 0.004621411907183619
$ws.Range("E8").Value = -0.00381643422540523
$ws.Range("D9").Value = 0.007007892990048887
$ws.Range("E9").Value = -0.03179929689996797
$ws.Range("D10").Value = 0.01452866788956444
$ws.Range("E10").Value = -0.01678356713426865
$ws.Range("D11").Value = 0.008516925449537416
$ws.Range("E11").Value = -0.01964937910883846
$ws.Range("D12").Value = 0.01576719931286167
$ws.Range("E12").Value = -0.01444128787878785
$ws.Range("D13").Value = 0.002860075640824389
$ws.Range("E13").Value = -0.01484583174724019
$ws.Range("D14").Value = 0.005854530821247031
$ws.Range("E14").Value = -0.007013442431326733
$ws.Range("D15").Value = 0.01454341231127036
$ws.Range("E15").Value = -0.01662324773601298
$ws.Range("D16").Value = 0.0105165987604341
$ws.Range("E16").Value = -0.003460668175162995
$ws.Range("D17").Value = 0.02058106636159492
$ws.Range("E17").Value = 0.005509079779636883
$ws.Range("D18").Value = 0.00845284631722899
$ws.Range("E18").Value = -0.006800618238021694
$ws.Range("D19").Value = 0.01694874385732927
$ws.Range("E19").Value = -0.00816350502143659
$ws.Range("D20").Value = 0.0121165551537758
$ws.Range("E20").Value = -0.00005134524543026764
$ws.Range("D21").Value = 0.007511007497794754
$ws.Range("E21").Value = -0.01133925835121075
$ws.Range("D22").Value = 0.01495470080029075
$ws.Range("E22").Value = -0.01723937099592299
$ws.Range("D23").Value = 0.02007798296022396
$ws.Range("E23").Value = -0.01125240525391114
$ws.Range("D24").Value = 0.01023325079140621
$ws.Range("E24").Value = -0.008021861777150918
$ws.Range("D25").Value = 0.02025796444543588
$ws.Range("E25").Value = -0.01189715745768116
$ws.Range("D26").Value = 0.01411834369816616
$ws.Range("E26").Value = -0.00832832091796587
$ws.Range("D27").Value = 0.02006656692063309
$ws.Range("E27").Value = 0.01019694773637969
$ws.Range("D28").Value = 0.05524181119734659
$ws.Range("E28").Value = -0.007410327158060781
$ws.Range("D29").Value = 0.0212495112519715
$ws.Range("E29").Value = -0.002210433244916099
$ws.Range("D30").Value = 0.02896625631340118
$ws.Range("E30").Value = 0.001932989690721643
$ws.Range("D31").Value = 0.01484016712788317
$ws.Range("E31").Value = -0.004420660772452178
$ws.Range("D32").Value = 0.01351074287710767
$ws.Range("E32").Value = -0.006561679790026198
$ws.Range("D33").Value = 0.01810764296086457
$ws.Range("E33").Value = -0.002102659245516225
$ws.Range("D34").Value = 0.04277282081587445
$ws.Range("E34").Value = -0.009464383113341301
$ws.Range("D35").Value = 0.01100419118710137
$ws.Range("E35").Value = -0.004748982360922582
$ws.Range("D36").Value = 0.01012251209673728
$ws.Range("E36").Value = -0.01144381345723622
$ws.Range("D37").Value = 0.01059871959018859
$ws.Range("E37").Value = -0.02003081664098605
$ws.Range("D38").Value = 0.007355413410509925
$ws.Range("E38").Value = 0.00786602385181423
$ws.Range("D39").Value = 0.01216461450300712
$ws.Range("E39").Value = -0.0114149821640902
$ws.Range("D40").Value = 0.01752831783459655
$ws.Range("E40").Value = -0.006207674943566666
$ws.Range("D41").Value = 0.01742566679740344
$ws.Range("E41").Value = -0.01126034459367786
$ws.Range("D42").Value = 0.03159972201302805
$ws.Range("E42").Value = 0.01086763924162781
$ws.Range("D43").Value = 0.01152292607207054
$ws.Range("E43").Value = -0.01060593004743382
$ws.Range("D44").Value = 0.02179023336697994
$ws.Range("E44").Value = -0.002168429437535924
$ws.Range("D45").Value = 0.0122089410872496
$ws.Range("E45").Value = 0.009103416647557871
$ws.Range("D46").Value = 0.008658241710739929
$ws.Range("E46").Value = -0.01173731690755653
$ws.Range("D47").Value = 0.01358409170913839
$ws.Range("E47").Value = -0.01364329582457435
$ws.Range("D48").Value = 0.01048832306564364
$ws.Range("E48").Value = 0.02818705957719425
$ws.Range("D49").Value = 0.01600329669840583
$ws.Range("E49").Value = -0.01556942095749025
$ws.Range("D50").Value = 0.008654477839376182
$ws.Range("E50").Value = -0.01643651472565077
$ws.Range("D51").Value = 0.01234643126433855
$ws.Range("E51").Value = -0.04643370033508865
$ws.Range("D52").Value = 0.008438941767645558
$ws.Range("E52").Value = -0.02201678627904147
$ws.Range("D53").Value = 0.01023334411053093
$ws.Range("E53").Value = -0.03023597250888366
$ws.Range("D54").Value = 0.1350524327002412
$ws.Range("E54").Value = 0.0001971608832807004
$ws.Range("D55").Value = 0.04389638307751539
$ws.Range("E55").Value = -0.007589447054571874
$ws.Range("E56").Value = -0.005326338976178846
